$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the row for "LINDA TANIA MIROSLAVA ENRIQUEZ PRECIADO" (original sheet row 10).
#    This shifts all subsequent rows up by one and drops the now-unused shared string.
$ws.Rows.Item(10).Delete()

# 2) The table shrank by one row when the sheet row was deleted; restore it to A1:J12
#    so row 12 exists again (blank, styled) exactly like the target workbook.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J12"))

# 3) Refresh the data for the remaining rows (2-11) with the updated figures.
#    Columns: A Asesor, B Clave, C Fecha_Corte, D Fecha_Conexion, E Mes_Asesor,
#             F Limite_Logro_Meta, G Polizas_Totales, H Comisones (I/J are formulas).

$ws.Range("C2").Value = 46059
$ws.Range("E2").Value = 13
$ws.Range("H2").Value = 99074.55

$ws.Range("C3").Formula = "=C2"
$ws.Range("E3").Value = 11
$ws.Range("H3").Value = 44177.39

$ws.Range("C4").Formula = "=C3"
$ws.Range("E4").Value = 10
$ws.Range("H4").Value = 96223.59

$ws.Range("C5").Formula = "=C4"
$ws.Range("E5").Value = 8
$ws.Range("H5").Value = 68894.89

$ws.Range("C6").Formula = "=C5"
$ws.Range("E6").Value = 8
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 93008.35

$ws.Range("C7").Formula = "=C6"
$ws.Range("E7").Value = 6
$ws.Range("H7").Value = 22109.76

$ws.Range("C8").Formula = "=C7"
$ws.Range("E8").Value = 6
$ws.Range("H8").Value = 30823.33

$ws.Range("C9").Formula = "=C8"
$ws.Range("E9").Value = 5
$ws.Range("G9").Value = 13.5
$ws.Range("H9").Value = 81420.55

# Row 10 now holds what used to be row 11 (ANA VERONICA GONZALEZ GAYTAN) - refresh
# its Clave/Fecha_Conexion/Mes_Asesor/Limite_Logro_Meta/Polizas_Totales/Comisones.
$ws.Range("B10").Value = 115404
$ws.Range("C10").Formula = "=C9"
$ws.Range("D10").Value = 45986
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 46296
$ws.Range("G10").Value = 9
$ws.Range("H10").Value = 48596.58

# Row 11 now holds what used to be row 12 (LUVIA PATRICIA FIGUEROA CASTRO) - refresh.
$ws.Range("B11").Value = 116060
$ws.Range("C11").Formula = "=C10"
$ws.Range("D11").Value = 46010
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 46327
$ws.Range("G11").Value = 7.5
$ws.Range("H11").Value = 36108.55

# 4) Update the window/selection state on the sheet view: clear the frozen
#    top-left cell + custom zoom, and move the active selection to K5.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
$ws.Range("K5").Select()

$wb.RecalculateFullRebuild()
